$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -970

$ws.Range("H32").Value = 5749.75
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3674

$ws.Range("H53").Value = 5442.125
$ws.Range("J53").Value = 546.75
$ws.Range("L53").Value = 546.75
$ws.Range("N53").Value = -1820.75

$ws.Range("H62").Value = 21298.166
$ws.Range("I62").Value = 3740
$ws.Range("J62").Value = 28051.309
$ws.Range("K62").Value = 3740
$ws.Range("L62").Value = 28051.309
$ws.Range("M62").Value = -3116
$ws.Range("N62").Value = -29299.309

$ws.Range("H65").Value = 21298.166
$ws.Range("I65").Value = 3740
$ws.Range("J65").Value = 28051.309
$ws.Range("K65").Value = 18700
$ws.Range("L65").Value = 140256.545
$ws.Range("M65").Value = -15580
$ws.Range("N65").Value = -146496.545

$ws.Range("H94").Value = 83543170
$ws.Range("I94").Value = 166676340
$ws.Range("K94").Value = 166676340
$ws.Range("M94").Value = -166675889

$ws.Range("H116").Value = 834992.75
$ws.Range("I116").Value = 3531880.2
$ws.Range("J116").Value = 5181.231
$ws.Range("K116").Value = 3531880.2
$ws.Range("L116").Value = 5181.231
$ws.Range("M116").Value = -3528438.2
$ws.Range("N116").Value = -12065.231

$ws.Range("H132").Value = 6346.4707
$ws.Range("I132").Value = 6332.6665
$ws.Range("K132").Value = 18997.9995
$ws.Range("M132").Value = -16467.9995

$ws.Range("H138").Value = 3647.276
$ws.Range("J138").Value = 4706.675
$ws.Range("L138").Value = 14120.025
$ws.Range("N138").Value = -24400.025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 160282.31
$ws.Range("I45").Value = 291423.56
$ws.Range("J45").Value = 7284.1665
$ws.Range("K45").Value = 291423.56
$ws.Range("L45").Value = 7284.1665
$ws.Range("M45").Value = -291046.56
$ws.Range("N45").Value = -8038.1665

$ws.Range("H61").Value = 5556.654
$ws.Range("I61").Value = 5905.5293
$ws.Range("J61").Value = 4897.6665
$ws.Range("K61").Value = 5905.5293
$ws.Range("L61").Value = 4897.6665
$ws.Range("M61").Value = -5693.5293
$ws.Range("N61").Value = -5321.6665

$ws.Range("H81").Value = 55000
$ws.Range("I81").Value = 55000
$ws.Range("K81").Value = 55000
$ws.Range("M81").Value = -54002

$ws.Range("H84").Value = 55000
$ws.Range("I84").Value = 55000
$ws.Range("K84").Value = 165000
$ws.Range("M84").Value = -160008

$ws.Range("H133").Value = 68999.5
$ws.Range("J133").Value = 68999.5
$ws.Range("L133").Value = 68999.5
$ws.Range("N133").Value = -74059.5

$ws.Range("H136").Value = 5556.654
$ws.Range("I136").Value = 5905.5293
$ws.Range("J136").Value = 4897.6665
$ws.Range("K136").Value = 17716.5879
$ws.Range("L136").Value = 14692.9995
$ws.Range("M136").Value = -15166.5879
$ws.Range("N136").Value = -19792.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 6790786
$ws.Range("I7").Value = 4255.2856
$ws.Range("K7").Value = 4255.2856
$ws.Range("M7").Value = -4142.2856

$ws.Range("H38").Value = 8433
$ws.Range("J38").Value = 8433
$ws.Range("L38").Value = 8433
$ws.Range("N38").Value = -9265

$ws.Range("H86").Value = 15938.556
$ws.Range("I86").Value = 14292.4
$ws.Range("J86").Value = 17996.25
$ws.Range("K86").Value = 14292.4
$ws.Range("L86").Value = 17996.25
$ws.Range("M86").Value = -13169.4
$ws.Range("N86").Value = -20242.25

$ws.Range("H89").Value = 15938.556
$ws.Range("I89").Value = 14292.4
$ws.Range("J89").Value = 17996.25
$ws.Range("K89").Value = 71462
$ws.Range("L89").Value = 89981.25
$ws.Range("M89").Value = -65846
$ws.Range("N89").Value = -101213.25

$ws.Range("H99").Value = 14865.129
$ws.Range("I99").Value = 20175.95
$ws.Range("K99").Value = 20175.95
$ws.Range("M99").Value = -18677.95

$ws.Range("H107").Value = 2621.3704
$ws.Range("I107").Value = 2684.3
$ws.Range("J107").Value = 2441.5715
$ws.Range("K107").Value = 2684.3
$ws.Range("L107").Value = 2441.5715
$ws.Range("M107").Value = -764.3000000000002
$ws.Range("N107").Value = -6281.5715

$ws.Range("H134").Value = 5943.3394
$ws.Range("I134").Value = 5989.26
$ws.Range("K134").Value = 17967.78
$ws.Range("M134").Value = -15432.78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4054.457
$ws.Range("I134").Value = 2217.2222
$ws.Range("J134").Value = 10255.125
$ws.Range("K134").Value = 6651.6666
$ws.Range("L134").Value = 30765.375
$ws.Range("M134").Value = -4116.6666
$ws.Range("N134").Value = -35835.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 3584.7
$ws.Range("I45").Value = 2949
$ws.Range("J45").Value = 3655.3333
$ws.Range("K45").Value = 8847
$ws.Range("L45").Value = 10965.9999
$ws.Range("M45").Value = -8315
$ws.Range("N45").Value = -12029.9999

$ws.Range("H62").Value = 3664.3333
$ws.Range("J62").Value = 8995
$ws.Range("L62").Value = 26985
$ws.Range("N62").Value = -28357

$ws.Range("H65").Value = 3664.3333
$ws.Range("J65").Value = 8995
$ws.Range("L65").Value = 80955
$ws.Range("N65").Value = -87819

$ws.Range("H137").Value = 2264.2104
$ws.Range("I137").Value = 2099.4375
$ws.Range("J137").Value = 3143
$ws.Range("K137").Value = 6298.3125
$ws.Range("L137").Value = 9429
$ws.Range("M137").Value = -1198.3125
$ws.Range("N137").Value = -19629

$ws.Range("H139").Value = 1768004.4
$ws.Range("I139").Value = 2310467.2
$ws.Range("K139").Value = 6931401.600000001
$ws.Range("M139").Value = -6926261.600000001

$ws.Range("H140").Value = 9344.637000000001
$ws.Range("I140").Value = 9646.762000000001
$ws.Range("K140").Value = 28940.286
$ws.Range("M140").Value = -23760.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 20004
$ws.Range("I14").Value = 20004
$ws.Range("K14").Value = 20004
$ws.Range("M14").Value = -19832

$ws.Range("H16").Value = 6740.609
$ws.Range("J16").Value = 344.33334
$ws.Range("L16").Value = 344.33334
$ws.Range("N16").Value = -684.33334

$ws.Range("H75").Value = 31999
$ws.Range("I75").Value = 31999
$ws.Range("K75").Value = 31999
$ws.Range("M75").Value = -31063

$ws.Range("H78").Value = 31999
$ws.Range("I78").Value = 31999
$ws.Range("K78").Value = 95997
$ws.Range("M78").Value = -91317

$ws.Range("H122").Value = 5988
$ws.Range("I122").Value = 5988
$ws.Range("K122").Value = 17964
$ws.Range("M122").Value = -15514

$ws.Range("H132").Value = 454254.9
$ws.Range("I132").Value = 710662.1
$ws.Range("J132").Value = 5542.3335
$ws.Range("K132").Value = 2131986.3
$ws.Range("L132").Value = 16627.0005
$ws.Range("M132").Value = -2129456.3
$ws.Range("N132").Value = -21687.0005

$ws.Range("H136").Value = 9155.923000000001
$ws.Range("I136").Value = 1487.25
$ws.Range("K136").Value = 4461.75
$ws.Range("M136").Value = -1911.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1720.4
$ws.Range("I17").Value = 1713
$ws.Range("K17").Value = 1713
$ws.Range("M17").Value = -1541

$ws.Range("H42").Value = 30044
$ws.Range("I42").Value = 30044
$ws.Range("K42").Value = 30044
$ws.Range("M42").Value = -29666

$ws.Range("H80").Value = 35650.5
$ws.Range("J80").Value = 35650.5
$ws.Range("L80").Value = 35650.5
$ws.Range("N80").Value = -37646.5

$ws.Range("H83").Value = 35650.5
$ws.Range("J83").Value = 35650.5
$ws.Range("L83").Value = 106951.5
$ws.Range("N83").Value = -116935.5
